# Finished the binary import/export functionality. Removed some example
# lines from the UDP configuration files. Import of dependent channel
# references now works.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample csv")

# Remove the three "FADEC 01 BITFIELD" example rows (FDC1BF00, FDC1BF01,
# FDC1BF02) that were only placeholders. Deleting the rows shifts the
# remaining signal rows up and shrinks the used range from A1:D49 to
# A1:D46.
$ws.Rows("5:7").Delete()

# Update the saved selection/active cell to reflect the new view state.
$ws.Range("C6").Select()
